$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New rows (7-17): persistence-control entries for the "Empresa" fields
# (Domicilio, Email, Empresa, Estado, Fiscales, Municipio) plus the
# existing module rows (Parametro, Rol, Subparametro, Usuario) re-listed.
$rows = @(
    @{ Row = 7;  B = "Domicilio";     Extra = $true;  Highlight = $true },
    @{ Row = 8;  B = "Email";         Extra = $false; Highlight = $true },
    @{ Row = 9;  B = "Empresa";       Extra = $false; Highlight = $false },
    @{ Row = 10; B = "Estado";        Extra = $false; Highlight = $false },
    @{ Row = 11; B = "Fiscales";      Extra = $false; Highlight = $true },
    @{ Row = 12; B = "Municipio";     Extra = $false; Highlight = $false },
    @{ Row = 13; B = "Parametro";     Extra = $false; Highlight = $false },
    @{ Row = 14; B = "Rol";           Extra = $false; Highlight = $false },
    @{ Row = 15; B = "Subparametro";  Extra = $false; Highlight = $false },
    @{ Row = 16; B = "Telefono";      Extra = $false; Highlight = $true },
    @{ Row = 17; B = "Usuario";       Extra = $false; Highlight = $false }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = "sistema"
    $ws.Range("B$n").Value = $r.B
    $ws.Range("E$n").Value = "si"
    if ($r.Extra) {
        $ws.Range("F$n").Value = "si"
    }
    if ($r.Highlight) {
        $ws.Range("B$n").Interior.Color = 65535
    }
}

# H8 stays empty but picks up a (no-visible-effect) cell format touch.
$ws.Range("H8").Interior.ColorIndex = -4142

[void]$ws.Range("B16").Select()

Write-Host "done"
